$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.084807892817217
$ws.Range("D2").Value = 1.087789444820161
$ws.Range("E2").Value = 1.081352253818684
$ws.Range("F2").Value = 1.096279409776523
$ws.Range("I2").Value = 1.070256735622021
$ws.Range("J2").Value = 1.089665855874333
$ws.Range("K2").Value = 1.090443803604722
$ws.Range("L2").Value = 1.084023300985792
$ws.Range("M2").Value = 1.098912081482161
$ws.Range("N2").Value = 1.091213306376095

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.086164336051116
$ws.Range("D3").Value = 1.088918187269163
$ws.Range("E3").Value = 1.082603774796071
$ws.Range("F3").Value = 1.097542191412052
$ws.Range("I3").Value = 1.070793534386877
$ws.Range("J3").Value = 1.090683394859031
$ws.Range("K3").Value = 1.091391432141117
$ws.Range("L3").Value = 1.085092202047608
$ws.Range("M3").Value = 1.09999500362776
$ws.Range("N3").Value = 1.092232290382864

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.087040873154767
$ws.Range("D4").Value = 1.089647392037441
$ws.Range("E4").Value = 1.083411358554033
$ws.Range("F4").Value = 1.098357844892488
$ws.Range("I4").Value = 1.071138959036448
$ws.Range("J4").Value = 1.091340106921377
$ws.Range("K4").Value = 1.092002852288965
$ws.Range("L4").Value = 1.085781055423364
$ws.Range("M4").Value = 1.100693709094737
$ws.Range("N4").Value = 1.092889935051665

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.087409094172977
$ws.Range("D5").Value = 1.089953674466145
$ws.Range("E5").Value = 1.083750337998497
$ws.Range("F5").Value = 1.098700403559619
$ws.Range("I5").Value = 1.071283718466109
$ws.Range("J5").Value = 1.091615785058391
$ws.Range("K5").Value = 1.092259476222516
$ws.Range("L5").Value = 1.086069985409672
$ws.Range("M5").Value = 1.100986966789208
$ws.Range("N5").Value = 1.093166004683255

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.087470904124491
$ws.Range("D6").Value = 1.090005084612483
$ws.Range("E6").Value = 1.083807223307106
$ws.Range("F6").Value = 1.098757900760593
$ws.Range("I6").Value = 1.071307997481592
$ws.Range("J6").Value = 1.091662049085772
$ws.Range("K6").Value = 1.092302540166764
$ws.Range("L6").Value = 1.08611845925388
$ws.Range("M6").Value = 1.101036178154176
$ws.Range("N6").Value = 1.093212334410861

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.087045794416488
$ws.Range("D7").Value = 1.089651485674671
$ws.Range("E7").Value = 1.083415890084128
$ws.Range("F7").Value = 1.098362423515703
$ws.Range("I7").Value = 1.071140895111913
$ws.Range("J7").Value = 1.091343792126788
$ws.Range("K7").Value = 1.092006282944492
$ws.Range("L7").Value = 1.08578491872192
$ws.Range("M7").Value = 1.100697629490831
$ws.Range("N7").Value = 1.09289362549049

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.085266555217197
$ws.Range("D8").Value = 1.08817115276699
$ws.Range("E8").Value = 1.081775677411837
$ws.Range("F8").Value = 1.096706476276319
$ws.Range("I8").Value = 1.070438548387036
$ws.Range("J8").Value = 1.090010093269207
$ws.Range("K8").Value = 1.090764425697315
$ws.Range("L8").Value = 1.084385123671483
$ws.Range("M8").Value = 1.099278480816785
$ws.Range("N8").Value = 1.091558032627554

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.082122082839736
$ws.Range("D9").Value = 1.085553492157601
$ws.Range("E9").Value = 1.078868031175178
$ws.Range("F9").Value = 1.09377713967603
$ws.Range("I9").Value = 1.069186102843675
$ws.Range("J9").Value = 1.087646705756896
$ws.Range("K9").Value = 1.088562464442718
$ws.Range("L9").Value = 1.081896829518227
$ws.Range("M9").Value = 1.096762085828633
$ws.Range("N9").Value = 1.089191288833892

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.080019218234254
$ws.Range("D10").Value = 1.083801986910249
$ws.Range("E10").Value = 1.076917531200989
$ws.Range("F10").Value = 1.091816288708044
$ws.Range("I10").Value = 1.068341011424757
$ws.Range("J10").Value = 1.086061942797654
$ws.Range("K10").Value = 1.087085060179769
$ws.Range("L10").Value = 1.080223035714655
$ws.Range("M10").Value = 1.095073639414363
$ws.Range("N10").Value = 1.087604275329473

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.079107023316265
$ws.Range("D11").Value = 1.083041988983911
$ws.Range("E11").Value = 1.076069994829036
$ws.Range("F11").Value = 1.090965256190445
$ws.Range("I11").Value = 1.067972638660366
$ws.Range("J11").Value = 1.085373487019988
$ws.Range("K11").Value = 1.086443034335227
$ws.Range("L11").Value = 1.079494643094662
$ws.Range("M11").Value = 1.094339880123851
$ws.Range("N11").Value = 1.086914841865636

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.078767940225135
$ws.Range("D12").Value = 1.082759448652222
$ws.Range("E12").Value = 1.075754730350092
$ws.Range("F12").Value = 1.090648842183334
$ws.Range("I12").Value = 1.06783543841205
$ws.Range("J12").Value = 1.085117421602423
$ws.Range("K12").Value = 1.086204206814835
$ws.Range("L12").Value = 1.079223533465914
$ws.Range("M12").Value = 1.094066924878339
$ws.Range("N12").Value = 1.086658412805807

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.078840686270367
$ws.Range("D13").Value = 1.082820065604477
$ws.Range("E13").Value = 1.075822376213516
$ws.Range("F13").Value = 1.090716727861831
$ws.Range("I13").Value = 1.06786488513554
$ws.Range("J13").Value = 1.085172364078642
$ws.Range("K13").Value = 1.086255452071483
$ws.Range("L13").Value = 1.079281712501425
$ws.Range("M13").Value = 1.094125493102484
$ws.Range("N13").Value = 1.086713433306644

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.079078999808524
$ws.Range("D14").Value = 1.083018639104638
$ws.Range("E14").Value = 1.076043944227519
$ws.Range("F14").Value = 1.0909391075214
$ws.Range("I14").Value = 1.0679613052158
$ws.Range("J14").Value = 1.085352327586115
$ws.Range("K14").Value = 1.086423299957638
$ws.Range("L14").Value = 1.079472244392998
$ws.Range("M14").Value = 1.094317325842739
$ws.Range("N14").Value = 1.086893652382939

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.079225798878294
$ws.Range("D15").Value = 1.083140954401048
$ws.Range("E15").Value = 1.076180399509975
$ws.Range("F15").Value = 1.091076082671401
$ws.Range("I15").Value = 1.068020663657216
$ws.Range("J15").Value = 1.085463163526628
$ws.Range("K15").Value = 1.08652666998105
$ws.Range("L15").Value = 1.079589564001108
$ws.Range("M15").Value = 1.094435466558138
$ws.Range("N15").Value = 1.087004645723201

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.080079722461124
$ws.Range("D16").Value = 1.083852391641464
$ws.Range("E16").Value = 1.076973716460887
$ws.Range("F16").Value = 1.091872726850583
$ws.Range("I16").Value = 1.068365407374044
$ws.Range("J16").Value = 1.086107585631395
$ws.Range("K16").Value = 1.087127620437044
$ws.Range("L16").Value = 1.080271299683087
$ws.Range("M16").Value = 1.095122280230615
$ws.Range("N16").Value = 1.087649982981272

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.080614922130593
$ws.Range("D17").Value = 1.08429822952927
$ws.Range("E17").Value = 1.077470546286531
$ws.Range("F17").Value = 1.092371908399711
$ws.Range("I17").Value = 1.068580999773102
$ws.Range("J17").Value = 1.086511210352064
$ws.Range("K17").Value = 1.087503961488479
$ws.Range("L17").Value = 1.080697958001149
$ws.Range("M17").Value = 1.0955523862904
$ws.Range("N17").Value = 1.088054180895341

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.080926937204209
$ws.Range("D18").Value = 1.084558126689424
$ws.Range("E18").Value = 1.077760053917595
$ws.Range("F18").Value = 1.092662882810613
$ws.Range("I18").Value = 1.068706515715782
$ws.Range("J18").Value = 1.086746421771088
$ws.Range("K18").Value = 1.087723253273999
$ws.Range("L18").Value = 1.080946470922037
$ws.Range("M18").Value = 1.095803004525329
$ws.Range("N18").Value = 1.08828972634156

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.081033299713377
$ws.Range("D19").Value = 1.084646719188062
$ws.Range("E19").Value = 1.077858720456267
$ws.Range("F19").Value = 1.092762065577552
$ws.Range("I19").Value = 1.068749273568311
$ws.Range("J19").Value = 1.086826586300207
$ws.Range("K19").Value = 1.087797988671422
$ws.Range("L19").Value = 1.081031148308245
$ws.Range("M19").Value = 1.095888415738657
$ws.Range("N19").Value = 1.088370004713505

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.080557516647904
$ws.Range("D20").Value = 1.0842504111598
$ws.Range("E20").Value = 1.077417270662202
$ws.Range("F20").Value = 1.092318370618105
$ws.Range("I20").Value = 1.068557893127534
$ws.Range("J20").Value = 1.086467927626372
$ws.Range("K20").Value = 1.087463606611906
$ws.Range("L20").Value = 1.080652217820266
$ws.Range("M20").Value = 1.095506266438691
$ws.Range("N20").Value = 1.088010836703214

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.079008829457972
$ws.Range("D21").Value = 1.082960170915322
$ws.Range("E21").Value = 1.075978710500194
$ws.Range("F21").Value = 1.090873630661049
$ws.Range("I21").Value = 1.067932922150792
$ws.Range("J21").Value = 1.085299342310294
$ws.Range("K21").Value = 1.086373882655327
$ws.Range("L21").Value = 1.079416152798941
$ws.Range("M21").Value = 1.094260847094063
$ws.Range("N21").Value = 1.086840591861949

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.078033639539889
$ws.Range("D22").Value = 1.08214753561862
$ws.Range("E22").Value = 1.075071613298359
$ws.Range("F22").Value = 1.089963511086738
$ws.Range("I22").Value = 1.067537834604158
$ws.Range("J22").Value = 1.084562623310369
$ws.Range("K22").Value = 1.085686699370663
$ws.Range("L22").Value = 1.078635791437681
$ws.Range("M22").Value = 1.093475459877992
$ws.Range("N22").Value = 1.086102826636541

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.078550747672071
$ws.Range("D23").Value = 1.082578464365138
$ws.Range("E23").Value = 1.075552733289303
$ws.Range("F23").Value = 1.090446151154787
$ws.Range("I23").Value = 1.067747482249199
$ws.Range("J23").Value = 1.084953361741608
$ws.Range("K23").Value = 1.086051182378694
$ws.Range("L23").Value = 1.079049781257632
$ws.Range("M23").Value = 1.093892032608375
$ws.Range("N23").Value = 1.086494119961171

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.080583456201761
$ws.Range("D24").Value = 1.084272018693624
$ws.Range("E24").Value = 1.07744134450304
$ws.Range("F24").Value = 1.09234256262292
$ws.Range("I24").Value = 1.068568334753343
$ws.Range("J24").Value = 1.086487485894222
$ws.Range("K24").Value = 1.087481841927882
$ws.Range("L24").Value = 1.080672886918989
$ws.Range("M24").Value = 1.095527106803026
$ws.Range("N24").Value = 1.088030422746048

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.082936133181687
$ws.Range("D25").Value = 1.086231328485736
$ws.Range("E25").Value = 1.079621826502635
$ws.Range("F25").Value = 1.094535821953365
$ws.Range("I25").Value = 1.069511663100859
$ws.Range("J25").Value = 1.08825929470065
$ws.Range("K25").Value = 1.089133367706296
$ws.Range("L25").Value = 1.082542718588488
$ws.Range("M25").Value = 1.097414523623965
$ws.Range("N25").Value = 1.089804747724211

